$d = $word.ActiveDocument

function Get-ParagraphByText($searchText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $rng.Paragraphs(1)
}

# --- Part 1: replace the empty paragraph after "JTAC Status" with the new
#     "CAP and SEAD" paragraph, and add a new bookmark-only paragraph after it.

# The _GoBack bookmark currently sits in the "================" paragraph
# right before the AWACS/DARKSTAR block; remove it from there first.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Locate the empty paragraph right after the "JTAC Status" bullet.
$jtacPara = Get-ParagraphByText("JTAC Status")
$emptyPara = $jtacPara.Next()
$emptyRange = $emptyPara.Range

# Make room for a second new paragraph (the bookmark-only one) right after it.
$emptyRange.InsertParagraphAfter()

# Refill the (still empty) paragraph with the new text (with proofErr spell-check markers).
$emptyPara = $jtacPara.Next()
$emptyRange = $emptyPara.Range
$emptyRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">-- CAP and SEAD support can be requested under the F10 / “Request </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Supoport</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” option.</w:t></w:r></w:p>')

# The newly made empty paragraph right after it gets the bookmark.
$bmPara = $emptyPara.Next()
$bmRange = $bmPara.Range
$bmRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

# --- Part 2: move <w:lastRenderedPageBreak/> from the "For example..." run
#     to the "All tanker freqs..." run.

$allTankerPara = Get-ParagraphByText("All tanker")
$allTankerPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">All tanker </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>freqs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> start with 317. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>( Flight</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> # 1-3 )( FL ##)</w:t></w:r></w:p>')

$forExamplePara = Get-ParagraphByText("For example")
$forExamplePara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">For example: TEXACO 3-1, FL </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>240  =</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">  317. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>( 3</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> )( 24 )</w:t></w:r></w:p>')
